{"js": "// Fill in the signature block at the end of the IP Agreement with the\n// signers' names (underlined) followed by a shortened underline/blank,\n// instead of a long blank line. Also restores the \"_GoBack\" bookmark\n// Word drops at the last edit location (right before the trailing\n// underscores on the Client's signature line).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the four \"Team Apple\" signer blank lines and the Client blank\n// line by their distinctive all-underscore text, scoped to the\n// ListParagraph-styled lines following the \"Team Apple\" / \"Client\"\n// headings (these are the only paragraphs whose entire text is a run of\n// underscores).\nconst blankLineIndexes = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (/^_+$/.test(t)) {\n    blankLineIndexes.push(i);\n  }\n}\n\nif (blankLineIndexes.length !== 5) {\n  throw new Error(\"Expected 5 blank signature lines, found \" + blankLineIndexes.length);\n}\n\nconst [teamIdx1, teamIdx2, teamIdx3, teamIdx4, clientIdx] = blankLineIndexes;\n\n// Helper: replace a blank \"_____\" paragraph with \"<Name><shorter blank>\",\n// underlining just the name while leaving the rest of the run properties\n// (size 36/half-point \"36pt\" = 18pt, etc.) intact by splitting the\n// existing formatted run via search() instead of inserting fresh runs.\nasync function fillSignatureLine(paragraph, name, trailingUnderscoreCount) {\n  const trailing = \"_\".repeat(trailingUnderscoreCount);\n  const whole = paragraph.getRange(\"Whole\");\n  whole.insertText(name + trailing, \"Replace\");\n  await context.sync();\n\n  const nameRanges = paragraph.search(name, { matchCase: true });\n  nameRanges.load(\"items\");\n  await context.sync();\n\n  nameRanges.items[0].font.underline = \"Single\";\n  await context.sync();\n}\n\n// Helper: same idea, but also drops a \"_GoBack\" bookmark partway through\n// the trailing underscore run (used for the very last edited line, which\n// is what Word stamps with _GoBack).\nasync function fillSignatureLineWithBookmark(paragraph, leading, name, mid, afterBookmark) {\n  // A \"|\" marker pinpoints the bookmark insertion point; it is removed\n  // again right before we drop the bookmark there.\n  const whole = paragraph.getRange(\"Whole\");\n  whole.insertText(leading + name + mid + \"|\" + afterBookmark, \"Replace\");\n  await context.sync();\n\n  const nameRanges = paragraph.search(name, { matchCase: true });\n  nameRanges.load(\"items\");\n  await context.sync();\n  nameRanges.items[0].font.underline = \"Single\";\n  await context.sync();\n\n  const markerRanges = paragraph.search(\"|\", { matchCase: true });\n  markerRanges.load(\"items\");\n  await context.sync();\n\n  const collapsed = markerRanges.items[0].insertText(\"\", \"Replace\");\n  await context.sync();\n  collapsed.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\nawait fillSignatureLine(paragraphs.items[teamIdx1], \"Christopher Jones\", 12);\nawait fillSignatureLine(paragraphs.items[teamIdx2], \"Dustin Cofer\", 16);\nawait fillSignatureLine(paragraphs.items[teamIdx3], \"Taylor Williams\", 13);\nawait fillSignatureLine(paragraphs.items[teamIdx4], \"Taisann Kham\", 15);\n\nawait fillSignatureLineWithBookmark(\n  paragraphs.items[clientIdx],\n  \"_\",\n  \"Evelyn R. Brannock\",\n  \"__\",\n  \"_______\"\n);\n", "ps1": "# Fill in the signature block at the end of the IP Agreement with the\n# signers' names (underlined) followed by a shortened underline/blank,\n# instead of a long blank line. Also restores the \"_GoBack\" bookmark\n# Word drops at the last edit location (right before the trailing\n# underscores on the Client's signature line).\n\n$d = $word.ActiveDocument\n\nfunction Repeat-Char($Char, $Count) {\n    $s = \"\"\n    for ($i = 0; $i -lt $Count; $i++) {\n        $s = $s + $Char\n    }\n    return $s\n}\n\nfunction Set-SignatureLine($Paragraph, $Name, $TrailingUnderscoreCount) {\n    $trailing = Repeat-Char \"_\" $TrailingUnderscoreCount\n    $Paragraph.Range.Text = $Name + $trailing\n\n    $searchRange = $Paragraph.Range.Duplicate\n    $find = $searchRange.Find\n    $find.ClearFormatting()\n    $find.Text = $Name\n    $find.MatchCase = $true\n    [void]$find.Execute()\n    $searchRange.Font.Underline = 1\n}\n\nfunction Set-SignatureLineWithBookmark($Paragraph, $Leading, $Name, $Mid, $AfterBookmark, $BookmarkName) {\n    $Paragraph.Range.Text = $Leading + $Name + $Mid + $AfterBookmark\n\n    $searchRange = $Paragraph.Range.Duplicate\n    $find = $searchRange.Find\n    $find.ClearFormatting()\n    $find.Text = $Name\n    $find.MatchCase = $true\n    [void]$find.Execute()\n    $searchRange.Font.Underline = 1\n\n    $nameEnd = $searchRange.End\n    $bookRange = $d.Range($nameEnd, $nameEnd + $Mid.Length)\n    $bookRange.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add($BookmarkName, $bookRange)\n}\n\n# Locate the five \"all underscores\" blank signature-line paragraphs (four\n# under \"Team Apple\", one under \"Client\") without hard-coding indexes.\n$blankParagraphs = @()\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $trimmed = $t.TrimEnd([char]13)\n    if ($trimmed -match '^_+$') {\n        $blankParagraphs += $p\n    }\n}\n\nif ($blankParagraphs.Count -ne 5) {\n    throw \"Expected 5 blank signature lines, found $($blankParagraphs.Count)\"\n}\n\nSet-SignatureLine $blankParagraphs[0] \"Christopher Jones\" 12\nSet-SignatureLine $blankParagraphs[1] \"Dustin Cofer\" 16\nSet-SignatureLine $blankParagraphs[2] \"Taylor Williams\" 13\nSet-SignatureLine $blankParagraphs[3] \"Taisann Kham\" 15\n\nSet-SignatureLineWithBookmark $blankParagraphs[4] \"_\" \"Evelyn R. Brannock\" \"__\" \"_______\" \"_GoBack\"\n"}
